$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.93 = 37151.79 pesos`n✅ 37151.79 pesos = 8.9 = 937.07 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- tasas: update the Binance/transfi rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 112
$ws2.Range("O10").Value = 4161
$ws2.Range("O12").Value = 105.3
